# Exponential Growth.xlsx edit
# - Add two more data rows (depth 9 & 10) to the "G0 - Dup Check" sheet
#   (Time and Position Count columns), which ripples through the
#   prediction/error formulas and the two scatter charts on that sheet.
# - Explicit-ize the shared formula in I14 on "G3 - None" (=6^B14).
# - Switch the active sheet / selections around (view-state only).

$wb = $excel.ActiveWorkbook

$wsNone = $wb.Worksheets.Item("G3 - None")
$wsG3Dup = $wb.Worksheets.Item("G3 - Dup Check")
$wsG0Dup = $wb.Worksheets.Item("G0 - Dup Check")

# --- "G0 - Dup Check": fill in the newly measured data points ---
$wsG0Dup.Range("C13").Value = 457.68599999999998
$wsG0Dup.Range("I13").Value = 4540668

$wsG0Dup.Range("C14").Value = 2411
$wsG0Dup.Range("I14").Value = 24563253

# The engine doesn't always re-type/re-evaluate formula cells that were
# previously cached as an empty string ("") once their precedents go from
# blank to populated, so nudge the dependent prediction/error formulas in
# rows 13 & 14 to force a fresh evaluation.
$wsG0Dup.Range("D13").Formula = $wsG0Dup.Range("D13").Formula
$wsG0Dup.Range("D14").Formula = $wsG0Dup.Range("D14").Formula
$wsG0Dup.Range("F13").Formula = $wsG0Dup.Range("F13").Formula
$wsG0Dup.Range("F14").Formula = $wsG0Dup.Range("F14").Formula
$wsG0Dup.Range("H13").Formula = $wsG0Dup.Range("H13").Formula
$wsG0Dup.Range("H14").Formula = $wsG0Dup.Range("H14").Formula

# --- "G3 - None": de-share the I14 formula (still =6^B14) ---
$wsNone.Range("I14").Formula = "=6^B14"

# --- View-state: selections on each sheet ---
$wsNone.Range("I17").Select()
$wsG3Dup.Range("I18").Select()
$wsG0Dup.Range("I21").Select()

# --- Active sheet switches from "G0 - Dup Check" to "G3 - Dup Check" ---
$wsG3Dup.Activate()
